$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.699.77"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.600.27"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.27%  "
Set-TextValue "D5" "211.43"
$ws.Range("E5").Value = "  +0.00%  "
Set-TextValue "D6" "0.513"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  +0.92%  "
Set-TextValue "D10" "19.54"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.824.86"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.596.66"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +0.44%  "
Set-TextValue "D16" "65.38"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "26.687.36"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "0.0₃0758"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "7.21"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.00"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "209.48"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue "D24" "8.95"
$ws.Range("E24").Value = "  +0.91%  "
Set-TextValue "D25" "142.85"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "1.291.84"
$ws.Range("E34").Value = "  +0.83%  "
Set-TextValue "D35" "0.619"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  +18.24%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  -1.01%  "
Set-TextValue "D42" "0.783"
$ws.Range("E42").Value = "  -0.02%  "
Set-TextValue "D43" "2.18"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "1.736.35"
Set-TextValue "D46" "90.99"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D49" "0.100"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0510"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  +0.15%  "
